$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B-column values per diff
$ws.Range("B7").Value = 1.724740706211402
$ws.Range("B8").Value = 0.9852216748768466
$ws.Range("B23").Value = 1.62833067638352
$ws.Range("B24").Value = 0.9747899159663916
$ws.Range("B30").Value = 0.203404346429716
$ws.Range("B31").Value = 0.1816239316239334
$ws.Range("B40").Value = -0.3816006600659966
$ws.Range("B41").Value = 0.6729475100942038
$ws.Range("B44").Value = 1.134157555941554
$ws.Range("B45").Value = 0.4041220448575383
$ws.Range("B46").Value = 0.2012477359629733
$ws.Range("B48").Value = 1.230123012301234
$ws.Range("B49").Value = 0.9583086346571812
$ws.Range("B52").Value = -0.4417554979352672
$ws.Range("B53").Value = 0.5112375807851849
$ws.Range("B54").Value = -0.5278310940499013
$ws.Range("B55").Value = 0.6753497346840218
$ws.Range("B56").Value = 0.287494010541458
$ws.Range("B57").Value = -0.1051122790253219
$ws.Range("B59").Value = -0.1524826074525977
$ws.Range("B60").Value = -1.737138493843651
$ws.Range("B68").Value = 0.5792972459639122
$ws.Range("B69").Value = -0.14162968558209
$ws.Range("B70").Value = 0.5673222390317647
$ws.Range("B72").Value = -0.3945885005636995
$ws.Range("B73").Value = -0.7074136955291455
$ws.Range("B74").Value = 0.1709888857224345
$ws.Range("B75").Value = 0.9009009009009035
$ws.Range("B76").Value = -0.7612781954887239
$ws.Range("B77").Value = -0.1420589070934801
$ws.Range("B78").Value = -0.3793626707131937
$ws.Range("B80").Value = 0.5155131264916527
$ws.Range("B81").Value = -0.09497578117580827

# Append new row 82 (copy formatting from row 81, then set values)
$ws.Range("A81:B81").Copy($ws.Range("A82"))
$ws.Range("A82").Value = 45884
$ws.Range("B82").Value = 0.06654624964350926
